$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update vm_pu values for rows 2-25 (bus indices 0-23), columns B-F and I-N
# matching the case with 380 kV re-run (commit: "case with 380 kV done")

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.055264595955309
$ws.Range("D2").Value = 1.055683548681641
$ws.Range("E2").Value = 1.063113272723415
$ws.Range("F2").Value = 1.071938796125796
$ws.Range("I2").Value = 1.048095089550095
$ws.Range("J2").Value = 1.060272364735297
$ws.Range("K2").Value = 1.058423331897787
$ws.Range("L2").Value = 1.065832791197818
$ws.Range("M2").Value = 1.074634637540979
$ws.Range("N2").Value = 1.061778073108161

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.056466819009624
$ws.Range("D3").Value = 1.05662280147425
$ws.Range("E3").Value = 1.064338039388658
$ws.Range("F3").Value = 1.073223380816733
$ws.Range("I3").Value = 1.04847273752535
$ws.Range("J3").Value = 1.061124975678202
$ws.Range("K3").Value = 1.059175753900822
$ws.Range("L3").Value = 1.066871476469536
$ws.Range("M3").Value = 1.075734717028411
$ws.Range("N3").Value = 1.062631894856401

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.057244191690191
$ws.Range("D4").Value = 1.057230030792895
$ws.Range("E4").Value = 1.065130946822328
$ws.Range("F4").Value = 1.074054927050838
$ws.Range("I4").Value = 1.048715650953296
$ws.Range("J4").Value = 1.061675585702415
$ws.Range("K4").Value = 1.059661466467438
$ws.Range("L4").Value = 1.067543407872147
$ws.Range("M4").Value = 1.076446316769421
$ws.Range("N4").Value = 1.063183286810004

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.057570870681871
$ws.Range("D5").Value = 1.057485184161318
$ws.Range("E5").Value = 1.065464382935604
$ws.Range("F5").Value = 1.074404590964825
$ws.Range("I5").Value = 1.048817425265145
$ws.Range("J5").Value = 1.061906803481513
$ws.Range("K5").Value = 1.059865384687535
$ws.Range("L5").Value = 1.067825849286299
$ws.Range("M5").Value = 1.076745421052916
$ws.Range("N5").Value = 1.06341483294487

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.057625714052685
$ws.Range("D6").Value = 1.057528018186095
$ws.Range("E6").Value = 1.065520374112555
$ws.Range("F6").Value = 1.074463305931342
$ws.Range("I6").Value = 1.048834493320429
$ws.Range("J6").Value = 1.061945610835258
$ws.Range("K6").Value = 1.059899607340072
$ws.Range("L6").Value = 1.067873270226694
$ws.Range("M6").Value = 1.076795638963167
$ws.Range("N6").Value = 1.063453695409508

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.057248557292921
$ws.Range("D7").Value = 1.05723344065783
$ws.Range("E7").Value = 1.065135401825456
$ws.Range("F7").Value = 1.07405959895306
$ws.Range("I7").Value = 1.048717012227297
$ws.Range("J7").Value = 1.061678676260392
$ws.Range("K7").Value = 1.059664192313197
$ws.Range("L7").Value = 1.067547182017693
$ws.Range("M7").Value = 1.076450313618696
$ws.Range("N7").Value = 1.063186381756928

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.055671007270701
$ws.Range("D8").Value = 1.056001083947161
$ws.Range("E8").Value = 1.063527106246172
$ws.Range("F8").Value = 1.072372859519089
$ws.Range("I8").Value = 1.048223018091756
$ws.Range("J8").Value = 1.060560734345998
$ws.Range("K8").Value = 1.058677856213756
$ws.Range("L8").Value = 1.066183855338001
$ws.Range("M8").Value = 1.075006462203061
$ws.Range("N8").Value = 1.062066852236785

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.052886880210853
$ws.Range("D9").Value = 1.053825406588282
$ws.Range("E9").Value = 1.060696082297388
$ws.Range("F9").Value = 1.069403081855965
$ws.Range("I9").Value = 1.047341401973651
$ws.Range("J9").Value = 1.058582396158091
$ws.Range("K9").Value = 1.056930911410704
$ws.Range("L9").Value = 1.063780148369736
$ws.Range("M9").Value = 1.072460411058442
$ws.Range("N9").Value = 1.060085704581734

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.051027773533126
$ws.Range("D10").Value = 1.052372123860624
$ws.Range("E10").Value = 1.058810639457981
$ws.Range("F10").Value = 1.067424770577315
$ws.Range("I10").Value = 1.046746128461133
$ws.Range("J10").Value = 1.057257779517197
$ws.Range("K10").Value = 1.055760228823582
$ws.Range("L10").Value = 1.06217667022167
$ws.Range("M10").Value = 1.070761725581489
$ws.Range("N10").Value = 1.058759206833274

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.050222010557074
$ws.Range("D11").Value = 1.051742148570754
$ws.Range("E11").Value = 1.057994649020165
$ws.Range("F11").Value = 1.066568477461776
$ws.Range("I11").Value = 1.04648657265621
$ws.Range("J11").Value = 1.056682827579849
$ws.Range("K11").Value = 1.055251857898702
$ws.Range("L11").Value = 1.061482085135707
$ws.Range("M11").Value = 1.070025838786242
$ws.Range("N11").Value = 1.058183438398227

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.049922597604771
$ws.Range("D12").Value = 1.051508041883669
$ws.Range("E12").Value = 1.05769161450686
$ws.Range("F12").Value = 1.066250458966814
$ws.Range("I12").Value = 1.046389890955547
$ws.Range("J12").Value = 1.056469055006838
$ws.Range("K12").Value = 1.055062805765875
$ws.Range("L12").Value = 1.061224042947681
$ws.Range("M12").Value = 1.069752443656294
$ws.Range("N12").Value = 1.057969362243642

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.049986828027263
$ws.Range("D13").Value = 1.051558263375375
$ws.Range("E13").Value = 1.057756613729435
$ws.Range("F13").Value = 1.066318672931528
$ws.Range("I13").Value = 1.046410641792025
$ws.Range("J13").Value = 1.0565149194744
$ws.Range("K13").Value = 1.055103368094676
$ws.Range("L13").Value = 1.061279395822459
$ws.Range("M13").Value = 1.069811090319294
$ws.Range("N13").Value = 1.058015291844009

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.0501972633866
$ws.Range("D14").Value = 1.051722799406304
$ws.Range("E14").Value = 1.057969598868175
$ws.Range("F14").Value = 1.066542189011576
$ws.Range("I14").Value = 1.046478586452739
$ws.Range("J14").Value = 1.056665161363426
$ws.Range("K14").Value = 1.055236235298239
$ws.Range("L14").Value = 1.06146075616383
$ws.Range("M14").Value = 1.070003240974081
$ws.Range("N14").Value = 1.05816574709375

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.050326904005168
$ws.Range("D15").Value = 1.051824161418628
$ws.Range("E15").Value = 1.058100834013295
$ws.Range("F15").Value = 1.066679910784796
$ws.Range("I15").Value = 1.046520413447279
$ws.Range("J15").Value = 1.05675770249029
$ws.Range("K15").Value = 1.055318069898111
$ws.Range("L15").Value = 1.061572492581898
$ws.Range("M15").Value = 1.070121624111368
$ws.Range("N15").Value = 1.058258419639629

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.051081233160956
$ws.Range("D16").Value = 1.052413918489717
$ws.Range("E16").Value = 1.058864802687026
$ws.Range("F16").Value = 1.067481606644474
$ws.Range("I16").Value = 1.046763316355268
$ws.Range("J16").Value = 1.057295907861401
$ws.Range("K16").Value = 1.055793936885303
$ws.Range("L16").Value = 1.062222761717097
$ws.Range("M16").Value = 1.070810556511954
$ws.Range("N16").Value = 1.0587973893241

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.051554199058673
$ws.Range("D17").Value = 1.052783670579036
$ws.Range("E17").Value = 1.059344130580011
$ws.Range("F17").Value = 1.067984575526828
$ws.Range("I17").Value = 1.046915200786779
$ws.Range("J17").Value = 1.057633137969497
$ws.Range("K17").Value = 1.056092044335387
$ws.Range("L17").Value = 1.062630585059636
$ws.Range("M17").Value = 1.071242611740161
$ws.Range("N17").Value = 1.05913509833763

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.05182999917285
$ws.Range("D18").Value = 1.052999273832547
$ws.Range("E18").Value = 1.059623754987163
$ws.Range("F18").Value = 1.068277980653182
$ws.Range("I18").Value = 1.047003618945013
$ws.Range("J18").Value = 1.057829705159011
$ws.Range("K18").Value = 1.056265784831801
$ws.Range("L18").Value = 1.062868435750532
$ws.Range("M18").Value = 1.07149458912993
$ws.Range("N18").Value = 1.0593319446751

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.051924027607804
$ws.Range("D19").Value = 1.05307277767058
$ws.Range("E19").Value = 1.059719106630384
$ws.Range("F19").Value = 1.068378029702388
$ws.Range("I19").Value = 1.047033737857675
$ws.Range("J19").Value = 1.057896706915262
$ws.Range("K19").Value = 1.056325002060136
$ws.Range("L19").Value = 1.062949532380155
$ws.Range("M19").Value = 1.071580501357628
$ws.Range("N19").Value = 1.05939904158153

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.051503461855713
$ws.Range("D20").Value = 1.052744006609144
$ws.Range("E20").Value = 1.059292699055805
$ws.Range("F20").Value = 1.067930608446824
$ws.Range("I20").Value = 1.046898922973752
$ws.Range("J20").Value = 1.057596970185812
$ws.Range("K20").Value = 1.056060074775496
$ws.Range("L20").Value = 1.06258683215525
$ws.Range("M20").Value = 1.07119625974542
$ws.Range("N20").Value = 1.059098879191544

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.050135298663225
$ws.Range("D21").Value = 1.051674350579299
$ws.Range("E21").Value = 1.057906878396677
$ws.Range("F21").Value = 1.066476367822517
$ws.Range("I21").Value = 1.046458585935235
$ws.Range("J21").Value = 1.05662092467896
$ws.Range("K21").Value = 1.055197115318711
$ws.Range("L21").Value = 1.061407351249243
$ws.Range("M21").Value = 1.06994665890416
$ws.Range("N21").Value = 1.058121447588118

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.049274403493517
$ws.Range("D22").Value = 1.051001202415065
$ws.Range("E22").Value = 1.057035907831849
$ws.Range("F22").Value = 1.065562299102899
$ws.Range("I22").Value = 1.0461801594988
$ws.Range("J22").Value = 1.056006031226071
$ws.Range("K22").Value = 1.054653261735415
$ws.Range("L22").Value = 1.060665518697709
$ws.Range("M22").Value = 1.069160671682197
$ws.Range("N22").Value = 1.057505680915996

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.049730845422302
$ws.Range("D23").Value = 1.051358109600801
$ws.Range("E23").Value = 1.057497593393822
$ws.Range("F23").Value = 1.066046839459527
$ws.Range("I23").Value = 1.046327907658462
$ws.Range("J23").Value = 1.056332113663102
$ws.Range("K23").Value = 1.054941690345303
$ws.Range("L23").Value = 1.061058802174035
$ws.Range("M23").Value = 1.069577368843057
$ws.Range("N23").Value = 1.057832226427491

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.051526388038971
$ws.Range("D24").Value = 1.052761929256557
$ws.Range("E24").Value = 1.059315938622476
$ws.Range("F24").Value = 1.067954993748794
$ws.Range("I24").Value = 1.046906278752771
$ws.Range("J24").Value = 1.05761331326272
$ws.Range("K24").Value = 1.05607452087776
$ws.Range("L24").Value = 1.062606602289058
$ws.Range("M24").Value = 1.071217204318989
$ws.Range("N24").Value = 1.059115245477495

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.053607165946388
$ws.Range("D25").Value = 1.05438836521949
$ws.Range("E25").Value = 1.061427625785764
$ws.Range("F25").Value = 1.070170560085693
$ws.Range("I25").Value = 1.047570644850865
$ws.Range("J25").Value = 1.059094846430591
$ws.Range("K25").Value = 1.057383600014091
$ws.Range("L25").Value = 1.064401734477293
$ws.Range("M25").Value = 1.073118850729386
$ws.Range("N25").Value = 1.060598882592399
